$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 11494831
$ws.Cells.Item(19, 9).Value = 33333604
$ws.Cells.Item(19, 11).Value = 33333604
$ws.Cells.Item(19, 13).Value = -33333429

$ws.Cells.Item(70, 8).Value = 1867.6154
$ws.Cells.Item(70, 9).Value = 1487.5
$ws.Cells.Item(70, 10).Value = 2036.5555
$ws.Cells.Item(70, 11).Value = 4462.5
$ws.Cells.Item(70, 12).Value = 6109.666499999999
$ws.Cells.Item(70, 13).Value = -4192.5
$ws.Cells.Item(70, 14).Value = -6649.666499999999

$ws.Cells.Item(73, 8).Value = 1867.6154
$ws.Cells.Item(73, 9).Value = 1487.5
$ws.Cells.Item(73, 10).Value = 2036.5555
$ws.Cells.Item(73, 11).Value = 4462.5
$ws.Cells.Item(73, 12).Value = 6109.666499999999
$ws.Cells.Item(73, 13).Value = -3526.5
$ws.Cells.Item(73, 14).Value = -7981.666499999999

$ws.Cells.Item(100, 8).Value = 1103.0555
$ws.Cells.Item(100, 9).Value = 959.6875
$ws.Cells.Item(100, 10).Value = 2250
$ws.Cells.Item(100, 11).Value = 959.6875
$ws.Cells.Item(100, 12).Value = 2250
$ws.Cells.Item(100, 13).Value = -418.6875
$ws.Cells.Item(100, 14).Value = -3332

$ws.Cells.Item(113, 8).Value = 2582.7778
$ws.Cells.Item(113, 9).Value = 2051.8667
$ws.Cells.Item(113, 10).Value = 5237.3335
$ws.Cells.Item(113, 11).Value = 2051.8667
$ws.Cells.Item(113, 12).Value = 5237.3335
$ws.Cells.Item(113, 13).Value = 1202.1333
$ws.Cells.Item(113, 14).Value = -11745.3335

$ws.Cells.Item(137, 8).Value = 569074.2
$ws.Cells.Item(137, 9).Value = 1988.9744
$ws.Cells.Item(137, 10).Value = 2780706.5
$ws.Cells.Item(137, 11).Value = 5966.9232
$ws.Cells.Item(137, 12).Value = 8342119.5
$ws.Cells.Item(137, 13).Value = -3416.9232
$ws.Cells.Item(137, 14).Value = -8347219.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 20162.844
$ws.Cells.Item(32, 9).Value = 21634.396
$ws.Cells.Item(32, 11).Value = 21634.396
$ws.Cells.Item(32, 13).Value = -21347.396

$ws.Cells.Item(135, 8).Value = 49428
$ws.Cells.Item(135, 10).Value = 49428
$ws.Cells.Item(135, 12).Value = 49428
$ws.Cells.Item(135, 14).Value = -59568

$ws.Cells.Item(139, 8).Value = 35854.668
$ws.Cells.Item(139, 10).Value = 35854.668
$ws.Cells.Item(139, 12).Value = 35854.668
$ws.Cells.Item(139, 14).Value = -46134.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1744.591
$ws.Cells.Item(86, 9).Value = 1801.55
$ws.Cells.Item(86, 10).Value = 1175
$ws.Cells.Item(86, 11).Value = 1801.55
$ws.Cells.Item(86, 12).Value = 1175
$ws.Cells.Item(86, 13).Value = -678.55
$ws.Cells.Item(86, 14).Value = -3421

$ws.Cells.Item(89, 8).Value = 1744.591
$ws.Cells.Item(89, 9).Value = 1801.55
$ws.Cells.Item(89, 10).Value = 1175
$ws.Cells.Item(89, 11).Value = 9007.75
$ws.Cells.Item(89, 12).Value = 5875
$ws.Cells.Item(89, 13).Value = -3391.75
$ws.Cells.Item(89, 14).Value = -17107

$ws.Cells.Item(107, 8).Value = 1399.3125
$ws.Cells.Item(107, 9).Value = 1005.5
$ws.Cells.Item(107, 11).Value = 1005.5
$ws.Cells.Item(107, 13).Value = 914.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 8465.643
$ws.Cells.Item(16, 9).Value = 12057.667
$ws.Cells.Item(16, 10).Value = 2000
$ws.Cells.Item(16, 11).Value = 12057.667
$ws.Cells.Item(16, 12).Value = 2000
$ws.Cells.Item(16, 13).Value = -11770.667
$ws.Cells.Item(16, 14).Value = -2574

$ws.Cells.Item(31, 8).Value = 4156.9375
$ws.Cells.Item(31, 9).Value = 4884.593
$ws.Cells.Item(31, 10).Value = 3221.3809
$ws.Cells.Item(31, 11).Value = 4884.593
$ws.Cells.Item(31, 12).Value = 3221.3809
$ws.Cells.Item(31, 13).Value = -4589.593
$ws.Cells.Item(31, 14).Value = -3811.3809

$ws.Cells.Item(34, 8).Value = 4156.9375
$ws.Cells.Item(34, 9).Value = 4884.593
$ws.Cells.Item(34, 10).Value = 3221.3809
$ws.Cells.Item(34, 11).Value = 4884.593
$ws.Cells.Item(34, 12).Value = 3221.3809
$ws.Cells.Item(34, 13).Value = -4682.593
$ws.Cells.Item(34, 14).Value = -3625.3809

$ws.Cells.Item(58, 8).Value = 2219336.8
$ws.Cells.Item(58, 9).Value = 3247899
$ws.Cells.Item(58, 10).Value = 3971.8462
$ws.Cells.Item(58, 11).Value = 3247899
$ws.Cells.Item(58, 12).Value = 3971.8462
$ws.Cells.Item(58, 13).Value = -3247696
$ws.Cells.Item(58, 14).Value = -4377.8462

$ws.Cells.Item(76, 8).Value = 9615.23
$ws.Cells.Item(76, 9).Value = 9615.23
$ws.Cells.Item(76, 11).Value = 9615.23
$ws.Cells.Item(76, 13).Value = -9300.23

$ws.Cells.Item(79, 8).Value = 9615.23
$ws.Cells.Item(79, 9).Value = 9615.23
$ws.Cells.Item(79, 11).Value = 9615.23
$ws.Cells.Item(79, 13).Value = -8523.23

$ws.Cells.Item(113, 8).Value = 8465.643
$ws.Cells.Item(113, 9).Value = 12057.667
$ws.Cells.Item(113, 10).Value = 2000
$ws.Cells.Item(113, 11).Value = 12057.667
$ws.Cells.Item(113, 12).Value = 2000
$ws.Cells.Item(113, 13).Value = -9887.666999999999
$ws.Cells.Item(113, 14).Value = -6340

$ws.Cells.Item(122, 8).Value = 12251.77
$ws.Cells.Item(122, 9).Value = 4938.25
$ws.Cells.Item(122, 10).Value = 100014
$ws.Cells.Item(122, 11).Value = 14814.75
$ws.Cells.Item(122, 12).Value = 300042
$ws.Cells.Item(122, 13).Value = -12364.75
$ws.Cells.Item(122, 14).Value = -304942

$ws.Cells.Item(134, 8).Value = 3422.7
$ws.Cells.Item(134, 9).Value = 2168.5
$ws.Cells.Item(134, 10).Value = 4258.8335
$ws.Cells.Item(134, 11).Value = 6505.5
$ws.Cells.Item(134, 12).Value = 12776.5005
$ws.Cells.Item(134, 13).Value = -3970.5
$ws.Cells.Item(134, 14).Value = -17846.5005

$ws.Cells.Item(136, 8).Value = 2219336.8
$ws.Cells.Item(136, 9).Value = 3247899
$ws.Cells.Item(136, 10).Value = 3971.8462
$ws.Cells.Item(136, 11).Value = 9743697
$ws.Cells.Item(136, 12).Value = 11915.5386
$ws.Cells.Item(136, 13).Value = -9741147
$ws.Cells.Item(136, 14).Value = -17015.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 5050902.5
$ws.Cells.Item(5, 9).Value = 388.03845
$ws.Cells.Item(5, 11).Value = 1164.11535
$ws.Cells.Item(5, 13).Value = -1052.11535

$ws.Cells.Item(12, 8).Value = 29411952
$ws.Cells.Item(12, 9).Value = 58823680
$ws.Cells.Item(12, 10).Value = 225.05882
$ws.Cells.Item(12, 11).Value = 176471040
$ws.Cells.Item(12, 12).Value = 675.17646
$ws.Cells.Item(12, 13).Value = -176470867
$ws.Cells.Item(12, 14).Value = -1021.17646

$ws.Cells.Item(33, 8).Value = 92.5
$ws.Cells.Item(33, 9).Value = 82.5
$ws.Cells.Item(33, 10).Value = 107.5
$ws.Cells.Item(33, 11).Value = 495
$ws.Cells.Item(33, 12).Value = 645
$ws.Cells.Item(33, 13).Value = -212
$ws.Cells.Item(33, 14).Value = -1211

$ws.Cells.Item(36, 8).Value = 1425.8
$ws.Cells.Item(36, 9).Value = 666.3333
$ws.Cells.Item(36, 10).Value = 2565
$ws.Cells.Item(36, 11).Value = 1998.9999
$ws.Cells.Item(36, 12).Value = 7695
$ws.Cells.Item(36, 13).Value = -1829.9999
$ws.Cells.Item(36, 14).Value = -8033

$ws.Cells.Item(39, 8).Value = 8409.666999999999
$ws.Cells.Item(39, 10).Value = 8409.666999999999
$ws.Cells.Item(39, 12).Value = 25229.001
$ws.Cells.Item(39, 14).Value = -25817.001

$ws.Cells.Item(41, 8).Value = 1043.8
$ws.Cells.Item(41, 9).Value = 79
$ws.Cells.Item(41, 10).Value = 1285
$ws.Cells.Item(41, 11).Value = 237
$ws.Cells.Item(41, 12).Value = 3855
$ws.Cells.Item(41, 13).Value = 101
$ws.Cells.Item(41, 14).Value = -4531

$ws.Cells.Item(80, 8).Value = 4500
$ws.Cells.Item(80, 10).Value = 5000
$ws.Cells.Item(80, 12).Value = 15000
$ws.Cells.Item(80, 14).Value = -16872

$ws.Cells.Item(83, 8).Value = 4500
$ws.Cells.Item(83, 10).Value = 5000
$ws.Cells.Item(83, 12).Value = 45000
$ws.Cells.Item(83, 14).Value = -54360

$ws.Cells.Item(97, 8).Value = 2203.111
$ws.Cells.Item(97, 10).Value = 2461.1428
$ws.Cells.Item(97, 12).Value = 7383.428400000001
$ws.Cells.Item(97, 14).Value = -8375.428400000001

$ws.Cells.Item(122, 8).Value = 867.1111
$ws.Cells.Item(122, 9).Value = 483.33334
$ws.Cells.Item(122, 10).Value = 894.5238000000001
$ws.Cells.Item(122, 11).Value = 4350.00006
$ws.Cells.Item(122, 12).Value = 8050.7142
$ws.Cells.Item(122, 13).Value = -1900.00006
$ws.Cells.Item(122, 14).Value = -12950.7142

$ws.Cells.Item(132, 8).Value = 2323.6924
$ws.Cells.Item(132, 9).Value = 2736.8572
$ws.Cells.Item(132, 10).Value = 1841.6666
$ws.Cells.Item(132, 11).Value = 24631.7148
$ws.Cells.Item(132, 12).Value = 16574.9994
$ws.Cells.Item(132, 13).Value = -22101.7148
$ws.Cells.Item(132, 14).Value = -21634.9994

$ws.Cells.Item(135, 8).Value = 5050902.5
$ws.Cells.Item(135, 9).Value = 388.03845
$ws.Cells.Item(135, 11).Value = 3492.34605
$ws.Cells.Item(135, 13).Value = -957.3460500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5337.54
$ws.Cells.Item(70, 9).Value = 4785.9565
$ws.Cells.Item(70, 11).Value = 4785.9565
$ws.Cells.Item(70, 13).Value = -4515.9565

$ws.Cells.Item(73, 8).Value = 5337.54
$ws.Cells.Item(73, 9).Value = 4785.9565
$ws.Cells.Item(73, 11).Value = 4785.9565
$ws.Cells.Item(73, 13).Value = -3849.9565

$ws.Cells.Item(80, 8).Value = 13777.777
$ws.Cells.Item(80, 9).Value = 37500
$ws.Cells.Item(80, 10).Value = 7000
$ws.Cells.Item(80, 11).Value = 37500
$ws.Cells.Item(80, 12).Value = 7000
$ws.Cells.Item(80, 13).Value = -36502
$ws.Cells.Item(80, 14).Value = -8996

$ws.Cells.Item(83, 8).Value = 13777.777
$ws.Cells.Item(83, 9).Value = 37500
$ws.Cells.Item(83, 10).Value = 7000
$ws.Cells.Item(83, 11).Value = 187500
$ws.Cells.Item(83, 12).Value = 35000
$ws.Cells.Item(83, 13).Value = -182508
$ws.Cells.Item(83, 14).Value = -44984

$ws.Cells.Item(113, 8).Value = 2016.36
$ws.Cells.Item(113, 9).Value = 1977.2142
$ws.Cells.Item(113, 10).Value = 2066.182
$ws.Cells.Item(113, 11).Value = 1977.2142
$ws.Cells.Item(113, 12).Value = 2066.182
$ws.Cells.Item(113, 13).Value = 192.7858000000001
$ws.Cells.Item(113, 14).Value = -6406.182

$ws.Cells.Item(132, 8).Value = 2498.65
$ws.Cells.Item(132, 9).Value = 2220.8333
$ws.Cells.Item(132, 11).Value = 6662.499899999999
$ws.Cells.Item(132, 13).Value = -4132.499899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 5903.1313
$ws.Cells.Item(122, 9).Value = 5132.2324
$ws.Cells.Item(122, 11).Value = 15396.6972
$ws.Cells.Item(122, 13).Value = -12946.6972

$ws.Cells.Item(132, 8).Value = 5778.7715
$ws.Cells.Item(132, 9).Value = 6382.16
$ws.Cells.Item(132, 10).Value = 4270.3
$ws.Cells.Item(132, 11).Value = 19146.48
$ws.Cells.Item(132, 12).Value = 12810.9
$ws.Cells.Item(132, 13).Value = -16616.48
$ws.Cells.Item(132, 14).Value = -17870.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 3871
$ws.Cells.Item(62, 10).Value = 3838.75
$ws.Cells.Item(62, 12).Value = 3838.75
$ws.Cells.Item(62, 14).Value = -5086.75

$ws.Cells.Item(65, 8).Value = 3871
$ws.Cells.Item(65, 10).Value = 3838.75
$ws.Cells.Item(65, 12).Value = 19193.75
$ws.Cells.Item(65, 14).Value = -25433.75

$ws.Cells.Item(122, 8).Value = 1780.4117
$ws.Cells.Item(122, 9).Value = 1417.4667
$ws.Cells.Item(122, 10).Value = 4502.5
$ws.Cells.Item(122, 11).Value = 4252.4001
$ws.Cells.Item(122, 12).Value = 13507.5
$ws.Cells.Item(122, 13).Value = -1802.4001
$ws.Cells.Item(122, 14).Value = -18407.5

$ws.Cells.Item(126, 8).Value = 1521.9445
$ws.Cells.Item(126, 9).Value = 1546
$ws.Cells.Item(126, 11).Value = 4638
$ws.Cells.Item(126, 13).Value = -2168

$ws.Cells.Item(132, 8).Value = 1893.1915
$ws.Cells.Item(132, 9).Value = 864.6177
$ws.Cells.Item(132, 11).Value = 2593.8531
$ws.Cells.Item(132, 13).Value = -63.85310000000027
